# add wait visibility to prevent missed collecting data
# Updates the "food" sheet: fills in missing Vietnamese descriptions,
# re-labels several Tea / Milk Tea / Blended Frappe rows, adjusts
# prices / categories / topping-category ids, and removes the two
# trailing rows (29-30) that are no longer part of the menu.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("food")

# --- Fill in previously-empty descriptions (column D) ---
$ws.Range("D4").Value  = "Latte Trà Xanh"
$ws.Range("D5").Value  = "Cà phê Almond Macchiato (Lớp Foam Phô Mai)"
$ws.Range("D9").Value  = "Cà phê Cappuccino"
$ws.Range("D10").Value = "Cà phê latte"

# --- Row 11: Mocha -> Strawberry & Pink Guava Tea ---
$ws.Range("B11").Value = "Strawberry & Pink Guava Tea"
$ws.Range("D11").Value = "Trà dâu ổi hồng"
$ws.Range("E11").Value = "TEA"
$ws.Range("F11").Value = 0

# --- Row 12: Caramel Caffé Macchiato (Cheese Foam) -> Peach Tea ---
$ws.Range("B12").Value = "Peach Tea"
$ws.Range("C12").Value = 38500
$ws.Range("D12").Value = "Trà Đào"
$ws.Range("E12").Value = "TEA"
$ws.Range("F12").Value = 0

# --- Row 13: Strawberry & Pink Guava Tea -> Fresh Fruit Tea ---
$ws.Range("B13").Value = "Fresh Fruit Tea"
$ws.Range("D13").Value = "Trà trái cây tươi"

# --- Row 14: Peach Tea -> Honey Tea ---
$ws.Range("B14").Value = "Honey Tea"
$ws.Range("D14").Value = "Trà mật ong"
$ws.Range("F14").Value = 2

# --- Row 15: Honey Tea -> Mango Tea ---
$ws.Range("B15").Value = "Mango Tea"
$ws.Range("D15").Value = "Trà Xoài"

# --- Row 16: Fresh Fruit Tea -> Oolong Milk Tea ---
$ws.Range("B16").Value = "Oolong Milk Tea"
$ws.Range("D16").Value = "Trà sữa Oolong"
$ws.Range("E16").Value = "MILK TEA"
$ws.Range("F16").Value = 3

# --- Row 17: Mango Tea -> Roasted Hazelnut Milk Tea ---
$ws.Range("B17").Value = "Roasted Hazelnut Milk Tea"
$ws.Range("D17").Value = "Trà sữa hạt phỉ"
$ws.Range("E17").Value = "MILK TEA"
$ws.Range("F17").Value = 3

# --- Row 18: Oolong Milk Tea -> Jasmine Milk Tea ---
$ws.Range("B18").Value = "Jasmine Milk Tea"
$ws.Range("D18").Value = "Trà sữa Lài nguyên lá"

# --- Row 19: Roasted Hazelnut Milk Tea -> Earl Grey Milk Tea ---
$ws.Range("B19").Value = "Earl Grey Milk Tea"
$ws.Range("D19").Value = "Trà Sữa Earl Grey"

# --- Row 20: Jasmine Milk Tea -> RaspBerry Milk Tea ---
$ws.Range("B20").Value = "RaspBerry Milk Tea"
$ws.Range("D20").Value = "Trà sữa dâu rừng"

# --- Row 21: Earl Grey Milk Tea -> Black Tea Full Leaf ---
$ws.Range("B21").Value = "Black Tea Full Leaf"
$ws.Range("D21").Value = "Trà sữa Trà Đen Nguyên Lá"
$ws.Range("F21").Value = 4

# --- Row 22: RaspBerry Milk Tea -> Macchiato with Cheese Foamed ---
$ws.Range("B22").Value = "Macchiato with Cheese Foamed"
$ws.Range("D22").Value = "Trà sữa macchiato"
$ws.Range("F22").Value = 5

# --- Row 23: Black Tea Full Leaf -> Cheesecake Frappé ---
$ws.Range("B23").Value = "Cheesecake Frappé"
$ws.Range("C23").Value = 41300
$ws.Range("D23").Value = "Bánh phô mai đá xay"
$ws.Range("E23").Value = "BLENDED FRAPPES"
$ws.Range("F23").Value = 6

# --- Row 24: Macchiato with Cheese Foamed -> Green Tea Frappé ---
$ws.Range("B24").Value = "Green Tea Frappé"
$ws.Range("C24").Value = 41300
$ws.Range("D24").Value = "Trà Xanh Đá Xay"
$ws.Range("E24").Value = "BLENDED FRAPPES"
$ws.Range("F24").Value = 6

# --- Row 25: Cheesecake Frappé -> Tiramisu Cake Frappé ---
$ws.Range("B25").Value = "Tiramisu Cake Frappé"
$ws.Range("D25").Value = "Bánh Tiramisu đá xay"

# --- Row 26: Green Tea Frappé -> Double Chocolate Frappé ---
$ws.Range("B26").Value = "Double Chocolate Frappé"
$ws.Range("D26").Value = "Sô-cô-la đá xay"

# --- Row 27: Tiramisu Cake Frappé -> Mocha Frappé ---
$ws.Range("B27").Value = "Mocha Frappé"
$ws.Range("D27").Value = "Cà Phê Mocha Đá Xay"

# --- Row 28: Double Chocolate Frappé -> Vanilla Caramel Frappé ---
$ws.Range("B28").Value = "Vanilla Caramel Frappé"
$ws.Range("D28").Value = "Vanilla caramel đá xay"

# --- Remove the two trailing rows that no longer belong to the menu ---
$ws.Range("A29:F30").EntireRow.Delete()
